# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Every player row gets the team's season record: 74 wins, 88 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1 ---
$headerValues = New-Object 'object[,]' 1,3
$headerValues[0,0] = "Wins"
$headerValues[0,1] = "Losses"
$headerValues[0,2] = "Ties"
$ws.Range("AD1:AF1").Value = $headerValues

# Match the formatting of the other header cells (bold, bordered, centered)
# by copying the style from the preceding header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-45): same record values repeated for each player ---
$recordValues = New-Object 'object[,]' 1,3
$recordValues[0,0] = 74
$recordValues[0,1] = 88
$recordValues[0,2] = 0

for ($r = 2; $r -le 45; $r++) {
    $ws.Range("AD$r`:AF$r").Value = $recordValues
}
